$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reformat the dates in column A from DD/MM/YYYY to DD-MM-YYYY and update
# the attendance counters (columns D, E, G, H) for the affected rows.
# A leading apostrophe is used for dates whose day component is <= 12 so
# that Excel keeps storing them as plain text instead of auto-converting
# them into date serial values.

$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("A4").Value = "'01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("A5").Value = "'04-08-2022"

$ws.Range("A6").Value = "'08-08-2022"

$ws.Range("A7").Value = "'11-08-2022"

$ws.Range("A8").Value = "15-08-2022"

$ws.Range("A9").Value = "18-08-2022"

$ws.Range("A10").Value = "22-08-2022"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

$ws.Range("A11").Value = "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

$ws.Range("A12").Value = "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

$ws.Range("A13").Value = "'01-09-2022"

$ws.Range("A14").Value = "'05-09-2022"

$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0

$ws.Range("A17").Value = "15-09-2022"

$ws.Range("A18").Value = "19-09-2022"

$ws.Range("A19").Value = "22-09-2022"

$ws.Range("A20").Value = "26-09-2022"

$ws.Range("A21").Value = "29-09-2022"

